$wb = $excel.ActiveWorkbook

# --- Restore the "tabSelected" view on Feuil3 to a plain scrolled view -----
# (Excel persists the current scroll position as topLeftCell on the sheet
#  that is no longer the active tab; we best-effort replicate that by
#  scrolling the window to row 3 before we move the active tab away.)
$ws3 = $wb.Worksheets.Item("Feuil3")
$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 3

# --- Add the new "Feuil4" worksheet at the end of the workbook ------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws4 = $wb.Worksheets.Add($null, $lastSheet)
$ws4.Name = "Feuil4"

# Row 1: hourly rate calc -> seconds
$ws4.Range("A1").Value = 0.023742
$ws4.Range("B1").Value = 110000
$ws4.Range("C1").Formula = "=B1*A1"
$ws4.Range("D1").Value = "seconds"

# Row 2: minutes
$ws4.Range("C2").Formula = "=C1/60"
$ws4.Range("D2").Value = "minutes"

# Row 3: hours
$ws4.Range("C3").Formula = "=C2/60"
$ws4.Range("D3").Value = "hours"

# Rows 4-7: weight fractions + their sum
$ws4.Range("H4").Value = 0.25
$ws4.Range("H5").Value = 0.333333
$ws4.Range("H6").Value = 0.103448
$ws4.Range("H7").Formula = "=SUM(H4:H6)"

# Row 8: misc checks
$ws4.Range("E8").Formula = "=1/1*(1/24)"
$ws4.Range("H8").Formula = "=H7/3"

# Row 9: misc check
$ws4.Range("C9").Formula = "=1/(26*4)"

# Selection on the new sheet matches the authored file (A2)
$ws4.Range("A2").Select()

# Make Feuil4 the active tab, as in the target workbook
$ws4.Activate()
